# Adds six new songs to the Songs-Index sheet, keeping the existing
# alphabetical-by-title ordering. Each entry is inserted as a new row
# (shifting subsequent rows down) directly before the row it now
# precedes, and the rank column (A) is rebuilt afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Songs-Index")

# (row to insert before -- in terms of the sheet as it is being edited,
#  i.e. already accounting for the earlier insertions above it --,
#  Song Title, Album, Artists, Genre, Year)
$newSongs = @(
    @(14, "Darts In The Dark", "Darts In The Dark (Single)", "MAGIC!", "Pop", 2017),
    @(40, "MIA", "MIA (Single)", "Bad Bunny, Drake", "Rap, hiphop", 2018),
    @(51, "Rewrite The Stars", "The Greatest Showman: Reimagined", "James Arthur, Anne-Marie", "Soundtrack", 2018),
    @(58, "Sweet But Psycho", "Sweet But Psycho (Single)", "Ava Max", "Pop", 2018),
    @(60, "Taste The Feeling", "Taste the Feeling - Single", "Avicii, Conrad Sewell", "Electronic", 2016),
    @(66, "Trampoline", "Trampoline", "SHAED", "Alternative", 2018)
)

foreach ($song in $newSongs) {
    $rowIndex = $song[0]

    $ws.Rows.Item($rowIndex).Insert()

    $ws.Cells.Item($rowIndex, 2).Value = $song[1]
    $ws.Cells.Item($rowIndex, 3).Value = $song[2]
    $ws.Cells.Item($rowIndex, 4).Value = $song[3]
    $ws.Cells.Item($rowIndex, 5).Value = $song[4]
    $ws.Cells.Item($rowIndex, 6).Value = $song[5]

    $ws.Cells.Item($rowIndex, 1).Style = $ws.Cells.Item($rowIndex - 1, 1).Style
}

# Renumber column A (the rank) for every data row now that the six new
# rows are in place.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
